$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new effort-log entry as row 45
$ws.Range("A45").Value = 41233
$ws.Range("A45").NumberFormat = "ddd\ dd/mm/yyyy"

$ws.Range("B45").Value = 2.5

$ws.Range("D45").Value = "New test case tc09 put to operation and completed. Fix of makefile. Documentation. Export for Andreas"

$ws.Range("D45").Select()
